# Update cryptos list (price / 1h volume) to the latest scraped snapshot.
# Some "Price" values look numeric (e.g. 575.76) but must stay text cells,
# exactly like the rest of column D; a leading "'" (Excel's quote-prefix)
# forces that while leaving the displayed text unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.746.53'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '3.073.57'
$ws.Range('E3').Value = '  -1.28%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''575.76'
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').Value = '''170.32'
$ws.Range('E6').Value = '  -0.65%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.069.24'
$ws.Range('E8').Value = '  -1.30%  '
$ws.Range('E9').Value = '  -1.99%  '
$ws.Range('E10').Value = '  -1.14%  '
$ws.Range('E11').Value = '  -1.58%  '
$ws.Range('E12').Value = '  -3.30%  '
$ws.Range('E13').Value = '  -2.38%  '
$ws.Range('D14').Value = '''35.70'
$ws.Range('E14').Value = '  -4.01%  '
$ws.Range('E15').Value = '  -1.90%  '
$ws.Range('D16').Value = '3.586.49'
$ws.Range('E16').Value = '  -1.14%  '
$ws.Range('D17').Value = '66.749.99'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = '''17.01'
$ws.Range('E18').Value = '  +3.97%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').Value = '''6.99'
$ws.Range('E19').Value = '  -2.43%  '
$ws.Range('D20').Value = '3.071.15'
$ws.Range('E20').Value = '  -1.38%  '
$ws.Range('D21').Value = '''489.38'
$ws.Range('E21').Value = '  +2.71%  '
$ws.Range('E22').Value = '  -3.01%  '
$ws.Range('D23').Value = '''0.688'
$ws.Range('E23').Value = '  -3.64%  '
$ws.Range('D24').Value = '''82.77'
$ws.Range('E24').Value = '  -1.54%  '
$ws.Range('D25').Value = '''12.68'
$ws.Range('E25').Value = '  -5.31%  '
$ws.Range('D26').Value = '''2.21'
$ws.Range('E26').Value = '  -3.30%  '
$ws.Range('D27').Value = '''10.17'
$ws.Range('E27').Value = '  +0.51%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('E29').Value = '  -0.82%  '
$ws.Range('D30').Value = '''2.27'
$ws.Range('E30').Value = '  -4.16%  '
$ws.Range('E31').Value = '  -2.44%  '
$ws.Range('D32').Value = '''27.54'
$ws.Range('E32').Value = '  -3.61%  '
$ws.Range('E33').Value = '  -2.69%  '
$ws.Range('E34').Value = '  -2.78%  '
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('B36').Value = 'Arweave'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D36').Value = '''47.43'
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('B37').Value = 'Mantle'
$ws.Range('C37').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D37').Value = '''0.948'
$ws.Range('E37').Value = '  -2.68%  '
$ws.Range('D38').Value = '''5.59'
$ws.Range('E38').Value = '  -4.57%  '
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('D40').Value = '''1.97'
$ws.Range('E40').Value = '  -4.74%  '
$ws.Range('D41').Value = '''0.300'
$ws.Range('E41').Value = '  -3.23%  '
$ws.Range('E42').Value = '  -4.78%  '
$ws.Range('D43').Value = '2.758.77'
$ws.Range('E43').Value = '  -2.52%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '''2.53'
$ws.Range('E44').Value = '  -2.13%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '''0.0345'
$ws.Range('E45').Value = '  -3.39%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = '''368.34'
$ws.Range('E46').Value = '  -4.34%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = '''135.22'
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('D49').Value = '''24.67'
$ws.Range('E49').Value = '  -0.71%  '
$ws.Range('E50').Value = '  -1.71%  '
$ws.Range('E51').Value = '  -2.02%  '
